# Fruta / hortaliza, semanal
# The dataset rows (2-41) get shuffled into a new row order while the
# header row (1) stays put. We snapshot the current values of every
# data row, then write them back out in the new order described by
# $rowMap (newRow -> oldRow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 41
$firstCol = 1   # A
$lastCol = 20   # T

# newRow -> oldRow : the value that used to live in oldRow now belongs in newRow
$rowMap = @{
    2 = 12
    3 = 24
    4 = 21
    5 = 11
    6 = 2
    7 = 31
    8 = 7
    9 = 19
    10 = 4
    11 = 25
    12 = 37
    13 = 41
    14 = 3
    15 = 38
    16 = 26
    17 = 35
    18 = 5
    19 = 10
    20 = 16
    21 = 23
    22 = 34
    23 = 39
    24 = 27
    25 = 15
    26 = 29
    27 = 9
    28 = 18
    29 = 22
    30 = 30
    31 = 28
    32 = 8
    33 = 20
    34 = 17
    35 = 36
    36 = 33
    37 = 13
    38 = 40
    39 = 32
    40 = 14
    41 = 6
}

# Snapshot every cell value in the data rows before we start overwriting.
# NOTE: use .Value() (explicit getter call) -- plain .Value resolves to a
# property-info object on this engine rather than invoking it.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

# Write back the snapshotted values in their new row positions.
for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $snapshot["$oldRow,$c"]
    }
}
